# Adds the "New Battery Storage" (New_BS) technology row to the ELC sheet
# of the SubRES_NEW_ELC template: fills in the previously-blank row 13
# (New Technologies block) and the previously-empty row 46 (Define
# Processes continuation) that references it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELC")

# ---------------------------------------------------------------------
# Row 13 - new "New Battery Storage" technology entry
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "New_BS"
$ws.Range("C13").Value = "New Battery Storage"
$ws.Range("D13").Value = "ELEC_HV"
$ws.Range("E13").Value = "ELEC_HV"
$ws.Range("F13").Value = 2025
$ws.Range("G13").Value = 0.8
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 2500
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = 25
$ws.Range("O13").Value = 1
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 1
$ws.Range("S13").Value = 1

# ---------------------------------------------------------------------
# Row 46 - continuation of the "Define Processes" table, pulling the new
# technology's short/long name in via formulas, plus its unit triple.
# ---------------------------------------------------------------------
$ws.Range("B46").Value = "STG"
$ws.Range("C46").Formula = "=B13"
$ws.Range("D46").Formula = "=C13"
$ws.Range("E46").Value = "PJ"
$ws.Range("F46").Value = "GW"
$ws.Range("G46").Value = "DAYNITE"

# Match the look of the sibling rows above (39-45): light-grey fill,
# small black Arial font.
$fmtRange = $ws.Range("B46,E46:G46")
$fmtRange.Interior.Color = 15921906
$fmtRange.Font.Name = "Arial"
$fmtRange.Font.Size = 10
$fmtRange.Font.Color = 0
